# "Actualizacion nov y correccion de archivos y permisos"
#
# Updates the reporting-period dates on row 8 of the "Reporte de
# Formatos" sheet (the visible/selected sheet) and moves the
# selection/scroll position back toward the start of the sheet
# (it had been left scrolled out at column T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8 / C8: "Fecha de inicio del periodo que se informa" / "Fecha de
# termino del periodo que se informa" -> updated to the new (Nov)
# reporting window.
$ws.Range("B8").Value = 43831   # 2020-01-01
$ws.Range("C8").Value = 44012   # 2020-06-30

# Scroll back toward the top-left of the sheet and leave the
# selection on D13 (previously scrolled/selected at T8).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 2
[void]$ws.Range("D13").Select()
